# Auto-generated script to update Hades_Profits leve-market data
# Applies per-cell numeric updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1237.8462
$ws.Range("I43").Value = 720
$ws.Range("J43").Value = 1361.1428
$ws.Range("K43").Value = 720
$ws.Range("L43").Value = 1361.1428
$ws.Range("M43").Value = -651
$ws.Range("N43").Value = -1499.1428
$ws.Range("H112").Value = 50002070
$ws.Range("J112").Value = 2343.8235
$ws.Range("L112").Value = 7031.470499999999
$ws.Range("N112").Value = -9247.470499999999
$ws.Range("H116").Value = 1850
$ws.Range("I116").Value = 1585.7142
$ws.Range("J116").Value = 2466.6667
$ws.Range("K116").Value = 1585.7142
$ws.Range("L116").Value = 2466.6667
$ws.Range("M116").Value = 1856.2858
$ws.Range("N116").Value = -9350.6667
$ws.Range("H127").Value = 1508.6364
$ws.Range("I127").Value = 448.75
$ws.Range("J127").Value = 2114.2856
$ws.Range("K127").Value = 1346.25
$ws.Range("L127").Value = 6342.8568
$ws.Range("M127").Value = 3613.75
$ws.Range("N127").Value = -16262.8568
$ws.Range("H137").Value = 1755855.6
$ws.Range("I137").Value = 2440036
$ws.Range("J137").Value = 2643.3125
$ws.Range("K137").Value = 7320108
$ws.Range("L137").Value = 7929.9375
$ws.Range("M137").Value = -7317558
$ws.Range("N137").Value = -13029.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2137.2285
$ws.Range("I2").Value = 1571
$ws.Range("J2").Value = 2986.5715
$ws.Range("K2").Value = 1571
$ws.Range("L2").Value = 2986.5715
$ws.Range("M2").Value = -1458
$ws.Range("N2").Value = -3212.5715
$ws.Range("H32").Value = 17078.152
$ws.Range("I32").Value = 16262.177
$ws.Range("K32").Value = 16262.177
$ws.Range("M32").Value = -15975.177
$ws.Range("H45").Value = 1156.6666
$ws.Range("I45").Value = 1045.8334
$ws.Range("K45").Value = 1045.8334
$ws.Range("M45").Value = -668.8334
$ws.Range("H61").Value = 17277208
$ws.Range("I61").Value = 20021082
$ws.Range("J61").Value = 127991
$ws.Range("K61").Value = 20021082
$ws.Range("L61").Value = 127991
$ws.Range("M61").Value = -20020870
$ws.Range("N61").Value = -128415
$ws.Range("H74").Value = 9316280
$ws.Range("I74").Value = 11145581
$ws.Range("J74").Value = 169775
$ws.Range("K74").Value = 11145581
$ws.Range("L74").Value = 169775
$ws.Range("M74").Value = -11144707
$ws.Range("N74").Value = -171523
$ws.Range("H77").Value = 9316280
$ws.Range("I77").Value = 11145581
$ws.Range("J77").Value = 169775
$ws.Range("K77").Value = 55727905
$ws.Range("L77").Value = 848875
$ws.Range("M77").Value = -55723537
$ws.Range("N77").Value = -857611
$ws.Range("H116").Value = 2137.2285
$ws.Range("I116").Value = 1571
$ws.Range("J116").Value = 2986.5715
$ws.Range("K116").Value = 1571
$ws.Range("L116").Value = 2986.5715
$ws.Range("M116").Value = 723
$ws.Range("N116").Value = -7574.5715
$ws.Range("H122").Value = 2585831.8
$ws.Range("I122").Value = 1899.4117
$ws.Range("K122").Value = 5698.2351
$ws.Range("M122").Value = -3248.2351
$ws.Range("H136").Value = 17277208
$ws.Range("I136").Value = 20021082
$ws.Range("J136").Value = 127991
$ws.Range("K136").Value = 60063246
$ws.Range("L136").Value = 383973
$ws.Range("M136").Value = -60060696
$ws.Range("N136").Value = -389073

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2137.2285
$ws.Range("I3").Value = 1571
$ws.Range("J3").Value = 2986.5715
$ws.Range("K3").Value = 1571
$ws.Range("L3").Value = 2986.5715
$ws.Range("M3").Value = -1457
$ws.Range("N3").Value = -3214.5715
$ws.Range("H107").Value = 3476.375
$ws.Range("I107").Value = 3587.2856
$ws.Range("K107").Value = 3587.2856
$ws.Range("M107").Value = -1667.2856
$ws.Range("H134").Value = 1632.0892
$ws.Range("I134").Value = 1047.762
$ws.Range("K134").Value = 3143.286
$ws.Range("M134").Value = -608.2860000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5956.25
$ws.Range("I31").Value = 2441.6667
$ws.Range("J31").Value = 16500
$ws.Range("K31").Value = 2441.6667
$ws.Range("L31").Value = 16500
$ws.Range("M31").Value = -2146.6667
$ws.Range("N31").Value = -17090
$ws.Range("H34").Value = 5956.25
$ws.Range("I34").Value = 2441.6667
$ws.Range("J34").Value = 16500
$ws.Range("K34").Value = 2441.6667
$ws.Range("L34").Value = 16500
$ws.Range("M34").Value = -2239.6667
$ws.Range("N34").Value = -16904
$ws.Range("H58").Value = 21278380
$ws.Range("I58").Value = 27779420
$ws.Range("J58").Value = 2246.182
$ws.Range("K58").Value = 27779420
$ws.Range("L58").Value = 2246.182
$ws.Range("M58").Value = -27779217
$ws.Range("N58").Value = -2652.182
$ws.Range("H107").Value = 411.77274
$ws.Range("I107").Value = 380.1875
$ws.Range("J107").Value = 496
$ws.Range("K107").Value = 380.1875
$ws.Range("L107").Value = 496
$ws.Range("M107").Value = 1539.8125
$ws.Range("N107").Value = -4336
$ws.Range("H132").Value = 16982.477
$ws.Range("I132").Value = 1418.4255
$ws.Range("J132").Value = 57621.945
$ws.Range("K132").Value = 4255.2765
$ws.Range("L132").Value = 172865.835
$ws.Range("M132").Value = -1725.2765
$ws.Range("N132").Value = -177925.835
$ws.Range("H134").Value = 17126.97
$ws.Range("I134").Value = 1202.6471
$ws.Range("J134").Value = 64899.94
$ws.Range("K134").Value = 3607.9413
$ws.Range("L134").Value = 194699.82
$ws.Range("M134").Value = -1072.9413
$ws.Range("N134").Value = -199769.82
$ws.Range("H136").Value = 21278380
$ws.Range("I136").Value = 27779420
$ws.Range("J136").Value = 2246.182
$ws.Range("K136").Value = 83338260
$ws.Range("L136").Value = 6738.545999999999
$ws.Range("M136").Value = -83335710
$ws.Range("N136").Value = -11838.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 755.65717
$ws.Range("I107").Value = 789.4286
$ws.Range("J107").Value = 733.1429000000001
$ws.Range("K107").Value = 2368.2858
$ws.Range("L107").Value = 2199.4287
$ws.Range("M107").Value = -448.2857999999997
$ws.Range("N107").Value = -6039.4287
$ws.Range("H123").Value = 2572.25
$ws.Range("I123").Value = 1911.25
$ws.Range("J123").Value = 3233.25
$ws.Range("K123").Value = 5733.75
$ws.Range("L123").Value = 9699.75
$ws.Range("M123").Value = -3283.75
$ws.Range("N123").Value = -14599.75
$ws.Range("H129").Value = 2875398.8
$ws.Range("I129").Value = 1289.9286
$ws.Range("K129").Value = 3869.7858
$ws.Range("M129").Value = 1130.2142
$ws.Range("H130").Value = 2749.0908
$ws.Range("I130").Value = 950
$ws.Range("K130").Value = 2850
$ws.Range("M130").Value = 2170
$ws.Range("H131").Value = 1116.3019
$ws.Range("J131").Value = 1220.3043
$ws.Range("L131").Value = 3660.9129
$ws.Range("N131").Value = -13740.9129
$ws.Range("H133").Value = 6071.3335
$ws.Range("I133").Value = 3501.111
$ws.Range("J133").Value = 7999
$ws.Range("K133").Value = 10503.333
$ws.Range("L133").Value = 23997
$ws.Range("M133").Value = -5443.332999999999
$ws.Range("N133").Value = -34117
$ws.Range("H134").Value = 5930.96
$ws.Range("I134").Value = 1935.5
$ws.Range("J134").Value = 7811.1763
$ws.Range("K134").Value = 5806.5
$ws.Range("L134").Value = 23433.5289
$ws.Range("M134").Value = -736.5
$ws.Range("N134").Value = -33573.5289
$ws.Range("H136").Value = 3142.8572
$ws.Range("H137").Value = 31277.65
$ws.Range("I137").Value = 788.5714
$ws.Range("J137").Value = 47694.848
$ws.Range("K137").Value = 2365.7142
$ws.Range("L137").Value = 143084.544
$ws.Range("M137").Value = 2734.2858
$ws.Range("N137").Value = -153284.544
$ws.Range("H138").Value = 3131.65
$ws.Range("I138").Value = 2045.4546
$ws.Range("J138").Value = 4459.222
$ws.Range("K138").Value = 6136.3638
$ws.Range("L138").Value = 13377.666
$ws.Range("M138").Value = -996.3638000000001
$ws.Range("N138").Value = -23657.666
$ws.Range("H139").Value = 4517.8037
$ws.Range("I139").Value = 2759
$ws.Range("J139").Value = 5058.974
$ws.Range("K139").Value = 8277
$ws.Range("L139").Value = 15176.922
$ws.Range("M139").Value = -3137
$ws.Range("N139").Value = -25456.922
$ws.Range("H140").Value = 2396.25
$ws.Range("I140").Value = 2396.25
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 7188.75
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -2008.75
$ws.Range("N140").Value = $null
$ws.Range("H141").Value = 7644.1816
$ws.Range("I141").Value = 7644.1816
$ws.Range("K141").Value = 22932.5448
$ws.Range("M141").Value = -17752.5448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1483.8462
$ws.Range("I113").Value = 996.6667
$ws.Range("J113").Value = 1630
$ws.Range("K113").Value = 996.6667
$ws.Range("L113").Value = 1630
$ws.Range("M113").Value = 1173.3333
$ws.Range("N113").Value = -5970

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 22111.307
$ws.Range("I132").Value = 1368.3793
$ws.Range("J132").Value = 52188.55
$ws.Range("K132").Value = 4105.1379
$ws.Range("L132").Value = 156565.65
$ws.Range("M132").Value = -1575.1379
$ws.Range("N132").Value = -161625.65
$ws.Range("H136").Value = 41153.86
$ws.Range("I136").Value = 26461.36
$ws.Range("J136").Value = 93245.45
$ws.Range("K136").Value = 79384.08
$ws.Range("L136").Value = 279736.35
$ws.Range("M136").Value = -76834.08
$ws.Range("N136").Value = -284836.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1136.579
$ws.Range("I126").Value = 1152.6471
$ws.Range("K126").Value = 3457.9413
$ws.Range("M126").Value = -987.9412999999995
$ws.Range("H132").Value = 42792.938
$ws.Range("I132").Value = 29351.371
$ws.Range("K132").Value = 88054.113
$ws.Range("M132").Value = -85524.113

Write-Host "Applied 252 cell updates."